$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.184479951858521
$ws.Range("B1").Value = 1.851878643035889
$ws.Range("C1").Value = 4.452043533325195
$ws.Range("D1").Value = 1.690649032592773
$ws.Range("E1").Value = 0.4554994702339172
